$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.215.20'
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.842.00'
$ws.Range("E3").Value = '  +0.28%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9993'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  +0.05%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6741'
$ws.Range("E6").Value = '  -1.70%  '
$ws.Range("E7").Value = '  +0.06%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07418'
$ws.Range("E8").Value = '  -0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2952'
$ws.Range("E9").Value = '  -1.94%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07720'
$ws.Range("E11").Value = '  +0.72%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.865.29'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("E13").Value = '  -1.10%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6702'
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '86.03'
$ws.Range("E15").Value = '  -1.80%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.115'
$ws.Range("E16").Value = '  -0.82%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.217.68'
$ws.Range("E17").Value = '  +0.29%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008306'
$ws.Range("E18").Value = '  +1.58%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.32'
$ws.Range("E19").Value = '  +0.20%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.50'
$ws.Range("E20").Value = '  -0.17%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.188'
$ws.Range("E22").Value = '  -2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.000'
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  +0.46%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '8.686'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1403'
$ws.Range("E26").Value = '  -3.72%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.00'
$ws.Range("E27").Value = '  -0.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.507'
$ws.Range("E28").Value = '  -0.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.174'
$ws.Range("E29").Value = '  -2.27%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.066'
$ws.Range("E30").Value = '  -2.11%  '
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.05310'
$ws.Range("E32").Value = '  +2.60%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7620'
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.871'
$ws.Range("E34").Value = '  +1.46%  '
$ws.Range("E35").Value = '  -0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.678'
$ws.Range("E36").Value = '  +0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.328.30'
$ws.Range("E37").Value = '  +1.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01804'
$ws.Range("E38").Value = '  -1.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.723'
$ws.Range("E39").Value = '  +0.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9189'
$ws.Range("E40").Value = '  -1.66%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.937'
$ws.Range("E41").Value = '  +2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '103.40'
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.08167'
$ws.Range("E44").Value = '  +16.33%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.993.85'
$ws.Range("B46").Value = 'BabyDogeCoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000124'
$ws.Range("E46").Value = '  +0.70%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5162'
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.775'
$ws.Range("E48").Value = '  +0.27%  '
$ws.Range("E49").Value = '  -2.02%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.161'
$ws.Range("E50").Value = '  -3.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05956'
$ws.Range("E51").Value = '  +0.33%  '
